$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the formula-driven logical test cells with their literal
# Boolean results (formulas removed, plain boolean constants written).

$ws.Range("A3").Value = $true
$ws.Range("B3").Value = $true
$ws.Range("C3").Value = $true
$ws.Range("D3").Value = $false
$ws.Range("E3").Value = $true

$ws.Range("A4").Value = $false
$ws.Range("B4").Value = $true
$ws.Range("C4").Value = $false
$ws.Range("D4").Value = $true
$ws.Range("E4").Value = $false

$ws.Range("A5").Value = $false
$ws.Range("B5").Value = $false

$ws.Range("A6").Value = $false
$ws.Range("B6").Value = $true

$ws.Range("A7").Value = $true
$ws.Range("B7").Value = $true

$ws.Range("A8").Value = $true
$ws.Range("B8").Value = $false
